# Resumo Casos de Uso - add new use case "CSU-04 / Editar Conta ADM"
# for the Administrador actor as a new row, right after the header block
# (new row 5), pushing every subsequent use-case row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 5; everything below shifts down.
$ws.Rows.Item(5).Insert()

# The inserted row comes back with default/no formatting, so clone the
# look of a normal data row (borders/alignment/etc., taken from what is
# now row 6) onto the new row 5.
$ws.Range("A6:D6").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new use case.
$ws.Cells.Item(5, 1).Value = "CSU-04"
$ws.Cells.Item(5, 2).Value = "Editar Conta ADM"
$ws.Cells.Item(5, 3).Value = "Administrador"

# Match the author's final selection state.
$ws.Range("D5").Select()
